# Update Financials for LYTS - Doing Updates for Financials
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LYTS")

# Row 48: Property Plant and Equipment
$ws.Range("D48").Value = 43700
$ws.Range("E48").Value = 47400
$ws.Range("F48").Value = 47500
$ws.Range("G48").Value = 43200
$ws.Range("H48").Value = 44300
$ws.Range("I48").Value = 45400
$ws.Range("J48").Value = 42500

# Row 49: Goodwill
$ws.Range("D49").Value = 65900
$ws.Range("E49").Value = 96700
$ws.Range("F49").Value = 16100
$ws.Range("G49").Value = 16600
$ws.Range("H49").Value = 17700
$ws.Range("I49").Value = 19100
$ws.Range("J49").Value = 24000
